# Applies the cryptos.xlsx data refresh described in the commit message:
# "Updated cryptos list on Sat May 25 07:33:14 UTC 2024 with GitHub Actions".
# Every D/E data cell in the sheet is stored as text (t="inlineStr" originally,
# no leading apostrophe shown to the user) even when its content happens to look
# like a plain number (e.g. "1.00", "6.35"). A handful of the new values below are
# exactly that kind of numeric-looking string, so for just those cells we first pin
# NumberFormat to Text ("@") to stop Excel's normal type-inference from silently
# turning "1.00" into the number 1 (losing the trailing zero) before writing them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textFormatRefs = @(
  "D4", "D5", "D6", "D11", "D14", "D21", "D22", "D25", "D27", "D28", "D30", "D32", "D33", "D36", "D38", "D42", "D43", "D44", "D48", "D49"
)
foreach ($ref in $textFormatRefs) {
  $ws.Range($ref).NumberFormat = "@"
}

# --- Cell value updates (row order matches the sheet) ---
$ws.Range("D2").Value = "68.782.15"
$ws.Range("E2").Value = "  +2.24%  "
$ws.Range("D3").Value = "3.755.45"
$ws.Range("E3").Value = "  +1.83%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "601.78"
$ws.Range("E5").Value = "  +1.52%  "
$ws.Range("D6").Value = "169.02"
$ws.Range("E6").Value = "  +1.31%  "
$ws.Range("D7").Value = "3.755.62"
$ws.Range("E7").Value = "  +1.90%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +2.64%  "
$ws.Range("E10").Value = "  +3.61%  "
$ws.Range("D11").Value = "6.35"
$ws.Range("E11").Value = "  +3.29%  "
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("D14").Value = "0.0000248"
$ws.Range("E14").Value = "  +3.22%  "
$ws.Range("D15").Value = "4.383.10"
$ws.Range("E15").Value = "  +1.99%  "
$ws.Range("D16").Value = "3.756.95"
$ws.Range("E16").Value = "  +2.08%  "
$ws.Range("D17").Value = "68.782.94"
$ws.Range("E18").Value = "  +2.69%  "
$ws.Range("E19").Value = "  +0.88%  "
$ws.Range("E20").Value = "  +1.82%  "
$ws.Range("D21").Value = "10.87"
$ws.Range("E21").Value = "  +20.22%  "
$ws.Range("D22").Value = "495.67"
$ws.Range("E22").Value = "  +2.41%  "
$ws.Range("E23").Value = "  +1.81%  "
$ws.Range("E24").Value = "  +9.63%  "
$ws.Range("D25").Value = "85.33"
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("E26").Value = "  +1.30%  "
$ws.Range("D27").Value = "12.36"
$ws.Range("E27").Value = "  +1.64%  "
$ws.Range("D28").Value = "10.29"
$ws.Range("E28").Value = "  +3.30%  "
$ws.Range("E29").Value = "  +0.38%  "
$ws.Range("D30").Value = "2.52"
$ws.Range("E30").Value = "  +7.45%  "
$ws.Range("E31").Value = "  +2.74%  "
$ws.Range("D32").Value = "7.91"
$ws.Range("E32").Value = "  +2.65%  "
$ws.Range("D33").Value = "31.93"
$ws.Range("E33").Value = "  +0.60%  "
$ws.Range("D34").Value = "3.900.93"
$ws.Range("D35").Value = "3.689.59"
$ws.Range("E35").Value = "  +1.80%  "
$ws.Range("D36").Value = "0.109"
$ws.Range("E36").Value = "  +1.90%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "1.01"
$ws.Range("E38").Value = "  +2.45%  "
$ws.Range("E39").Value = "  +1.92%  "
$ws.Range("E40").Value = "  +1.40%  "
$ws.Range("E41").Value = "  +0.92%  "
$ws.Range("D42").Value = "439.85"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").Value = "48.88"
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "2.93"
$ws.Range("E44").Value = "  +5.84%  "
$ws.Range("E45").Value = "  +2.85%  "
$ws.Range("E46").Value = "  +2.27%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "40.53"
$ws.Range("E48").Value = "  +2.38%  "
$ws.Range("D49").Value = "141.59"
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").Value = "2.789.71"
$ws.Range("E50").Value = "  +1.27%  "
$ws.Range("E51").Value = "  +2.81%  "
